$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "CITY" column header (E1), matching the style of the existing
#     header row (bold / bordered / centered) by copying D1's format ---
$ws.Cells.Item(1, 5).Value = "CITY"
$ws.Cells.Item(1, 4).Copy()
$ws.Cells.Item(1, 5).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Fill E2:E9 with blank values for the already-existing rows ---
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 5).Value = ""
}

# --- Append new row 10 of data ---
# Column A holds date-like text ("2025-09-13"); force it to stay plain
# text (matching every other cell in that column, which is text too)
# instead of letting Excel reinterpret it as a date serial number. The
# NumberFormat is reset back to the same (default) look as the sibling
# row right afterwards via a format-only paste.
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "2025-09-13"
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B10").Value = "AA"
$ws.Range("C10").Value = "44CDX12"
$ws.Range("D10").Value = "MAM "
$ws.Range("E10").Value = "Chennai"
